# Generate Report for Archive
#
# Update the localization status report: the two in-flight files
# (42f4766d-43ab-4701-a29c-f352d6e86e17.md and
#  b5df1cc1-dfc1-4509-9a88-1b812e811937.md) have moved from
# "Ready for handoff" to "In Translation" for both the zh-cn and de-de
# locales. Reflect this on the per-locale sheets as well as the Overview
# rollup sheet.

$wb = $excel.ActiveWorkbook

$newStatus = "In Translation"

# --- Overview sheet: columns B (zh-cn) and C (de-de), rows 3 & 4 ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = $newStatus
$overview.Range("C3").Value = $newStatus
$overview.Range("B4").Value = $newStatus
$overview.Range("C4").Value = $newStatus

# --- zh-cn sheet: Status column (C), rows 3 & 4 ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = $newStatus
$zhcn.Range("C4").Value = $newStatus

# --- de-de sheet: Status column (C), rows 3 & 4 ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = $newStatus
$dede.Range("C4").Value = $newStatus
